$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data runs through row 357 (date serial 44431, 2021-08-23).
# Append rows 358-366 continuing the daily series through 2021-09-01
# (serial 44440), matching the commit "aggiornamento fino a 1/09/2021".
$lastRow = 357
$startSerial = 44432
$newRowCount = 9

for ($i = 0; $i -lt $newRowCount; $i++) {
    $row = $lastRow + 1 + $i
    $serial = $startSerial + $i

    # Copy the formatting (style incl. number format/border/alignment) of
    # column A from the previous row so the new date cell matches s="2".
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

$excel.CutCopyMode = $false
